# Apply "PCB done, BOM not ready" update to the journal worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Extend the blank template rows further down the sheet -------------
# Rows 32:33 already hold the blank alternating row-style pattern
# (s8/s14/s8/s6 then s4/s4/s4/s5). Copy that two-row block down so the
# same pattern continues through row 39.
$ws.Range("B32:E33").Copy() | Out-Null
$ws.Range("B34:E39").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- 2. Give C31 the same "date" formatting already used by C15 -----------
# (row 31 was a blank "General" row; it now needs the date number format
# used on the other date cells in the s4-pattern rows).
$ws.Range("C15").Copy() | Out-Null
$ws.Range("C31").PasteSpecial(-4122) | Out-Null        # xlPasteFormats
$excel.CutCopyMode = $false

# --- 3. Fill in the new journal entries ------------------------------------
$ws.Range("B31").Value = "Samedi"
$ws.Range("C31").Value = 44996
$ws.Range("D31").Value = 4
$ws.Range("E31").Value = "Routage du PCB"

$ws.Range("B32").Value = "Mercredi"
$ws.Range("C32").Value = 45000
$ws.Range("D32").Value = 8
$ws.Range("E32").Value = "Routage du PCB, contrôle du PCB par Ali Zoubir et modification de la BOM"

# Match the author's final selection/view position in the saved file.
$ws.Range("E30").Select() | Out-Null

$wb.Save()
